# updated BOM with proper connector
# J2/J3 connector part: manufacturer "Phoenix Contact" -> "Molex",
# and MPN changed from the numeric KiCad footprint suffix (1725672)
# to the real Molex part number "70553-0038".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "Molex"
$ws.Range("G3").Value = "70553-0038"

# Widen the Reference column (A) so the longer part info is readable.
$ws.Columns.Item(1).ColumnWidth = 44.1666666666667

# Tweak the page top/bottom margins slightly.
$ws.PageSetup.TopMargin = 73.8
$ws.PageSetup.BottomMargin = 73.8

# Simplify the header/footer (drop the explicit Times New Roman styling).
$ws.PageSetup.CenterHeader = "&A"
$ws.PageSetup.CenterFooter = "Page &P"

# Leave the cursor parked on I13, matching where editing left off.
$ws.Range("I13").Select()
